# modul receiving update 10/24/25
#
# The inventory table's first data row (row 2, just under the header)
# through row 7 get straightforward Part Number replacements. Starting
# with row 8, several obsolete line items are removed outright (their
# whole row disappears) while the row that follows absorbs the new
# part number in its first cell, keeping its own qty/stock/location/name
# values untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1) Remove the rows that are dropped entirely from the report ---
# Delete from the bottom up so earlier row indexes stay valid.
$rowsToDelete = @(19, 18, 17, 14, 13, 12, 10, 8)
foreach ($idx in $rowsToDelete) {
    $t.Rows.Item($idx).Delete()
}

# --- 2) Straight Part Number text swaps (Find/Replace, whole word only) ---
$replacements = @(
    @("242193209", "WE12X27300"),
    @("WB06X10610", "D517191P"),
    @("WB03X10348", "6600JB1010A"),
    @("WB24X10146", "316238201"),
    @("WB24X10205", "WH12X22744"),
    @("WB24X10204", "WH44X10288"),
    @("WPW10535778", "WE4M416"),
    @("5304500204", "205765"),
    @("5304525218", "38174"),
    @("5303931775", "12112425")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
